# "updated mouseHover action keyword"
#
# Adds a second full pass of the Login test-step sequence ("Login_02") to
# the "Test Steps" sheet -- this run ends up hovering/clicking the Account
# menu and is marked FAIL -- and reflects that outcome back on the
# "Test Cases" sheet (Login_01 is flipped to Runmode "No" since it already
# ran; Login_02 is flipped to "Yes" to run next, and its prior Results cell
# now shows "FAIL").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Test Cases"
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Range("B2").Value = "Login in and log out into Amazon"
$wsCases.Range("C2").Value = "No"
$wsCases.Range("C3").Value = "Yes"
$wsCases.Range("D3").Value = "FAIL"

# ---------------------------------------------------------------------------
# "Test Steps" - append the Login_02 block (rows 13-21), mirroring the
# Login_01 block (rows 2-10) against the "zeba" data set.
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Test Steps")

# Each entry: row, B (TS_ID), C (Description), E (Page Object), F (Action
# Keyword), G (Data Set, optional), H (Results)
$steps = @(
    @(13, "TS_01", "Open the Browser", "", "openBrowser", "Chrome", "PASS"),
    @(14, "TS_02", "Navigate to website", "", "navigate", "", "PASS"),
    @(15, "TS_03", "Click on SignIn button ", "btn_SignIn", "click", "", "PASS"),
    @(16, "TS_04", "Enter the Email in the Email address field", "txtbx_Email", "input", "zebatanveer.2013@gmail.com", "PASS"),
    @(17, "TS_05", "Click on Continue button", "btn_Continue", "click", "", "PASS"),
    @(18, "TS_06", "Enter the Password in the Password field", "txtbx_Password", "input", "zeba", "PASS"),
    @(19, "TS_07", "Click on Login button", "btn_LogIn", "click", "", "PASS"),
    @(20, "TS_08", "wait for some time", "", "wait_For", "", "PASS"),
    @(21, "TS_09", "Click on Account and list button", "btn_Account", "mouseHover", "", "FAIL")
)

foreach ($step in $steps) {
    $row = $step[0]
    $wsSteps.Range("A$row").Value = "Login_02"
    $wsSteps.Range("B$row").Value = $step[1]
    $wsSteps.Range("C$row").Value = $step[2]
    if ($step[3] -ne "") { $wsSteps.Range("E$row").Value = $step[3] }
    $wsSteps.Range("F$row").Value = $step[4]
    if ($step[5] -ne "") { $wsSteps.Range("G$row").Value = $step[5] }
    $wsSteps.Range("H$row").Value = $step[6]
}

# Row 16's Data Set cell is the same e-mail address as row 5, and carries
# the same hyperlink + Hyperlink-style formatting.
$wsSteps.Range("G16").Style = $wsSteps.Range("G5").Style
$wsSteps.Hyperlinks.Add($wsSteps.Range("G16"), "mailto:zebatanveer.2013@gmail.com") | Out-Null
$wsSteps.Range("G16").Style = $wsSteps.Range("G5").Style

# Extend the "Page_Name" list validation down over the newly added rows.
$wsSteps.Range("F13:F21").Validation.Add(3, 1, 3, "=Page_Name") | Out-Null
$wsSteps.Range("D13:D21").Validation.Add(3, 1, 3, "=Page_Name") | Out-Null

# Keep "Test Steps" the active sheet, cursor on the last edited cell.
$wsSteps.Activate() | Out-Null
$wsSteps.Range("H21").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore the "Test Cases" selection without changing which sheet is active.
# ---------------------------------------------------------------------------
$wsCases.Range("C3").Select() | Out-Null
$wsSteps.Activate() | Out-Null
